$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D from Excel auto-converting numeric-looking strings to numbers
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Range("D2").Value = "43.991.92"
$ws.Range("E2").Value = "  +0.20%  "

$ws.Range("D3").Value = "2.235.76"
$ws.Range("E3").Value = "  -0.52%  "

$ws.Range("D4").Value = "1.01"
$ws.Range("E4").Value = "  +0.50%  "

$ws.Range("D5").Value = "305.33"
$ws.Range("E5").Value = "  -4.17%  "

$ws.Range("D6").Value = "94.79"
$ws.Range("E6").Value = "  -5.64%  "

$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  -0.67%  "

$ws.Range("E8").Value = "  +0.31%  "

$ws.Range("D9").Value = "0.520"
$ws.Range("E9").Value = "  -4.91%  "

$ws.Range("D10").Value = "34.79"
$ws.Range("E10").Value = "  -5.70%  "

$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -2.93%  "

$ws.Range("D12").Value = "7.20"
$ws.Range("E12").Value = "  -4.46%  "

$ws.Range("E13").Value = "  -1.43%  "

$ws.Range("D14").Value = "2.577.60"
$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("D15").Value = "2.237.28"
$ws.Range("E15").Value = "  -0.75%  "

$ws.Range("D16").Value = "0.825"
$ws.Range("E16").Value = "  -3.25%  "

$ws.Range("D17").Value = "13.63"
$ws.Range("E17").Value = "  -5.92%  "

$ws.Range("D18").Value = "43.898.43"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("D19").Value = "0.0₃0958"
$ws.Range("E19").Value = "  -2.08%  "

$ws.Range("D20").Value = "12.12"
$ws.Range("E20").Value = "  -10.04%  "

$ws.Range("E21").Value = "  -3.04%  "

$ws.Range("D22").Value = "64.86"
$ws.Range("E22").Value = "  -0.70%  "

$ws.Range("D23").Value = "236.22"
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("E24").Value = "  -5.62%  "

$ws.Range("D25").Value = "1.95"
$ws.Range("E25").Value = "  -5.23%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").Value = "10.01"
$ws.Range("E27").Value = "  -6.54%  "

$ws.Range("D28").Value = "37.41"
$ws.Range("E28").Value = "  -3.21%  "

$ws.Range("E29").Value = "  -1.53%  "

$ws.Range("D30").Value = "5.93"
$ws.Range("E30").Value = "  -3.13%  "

$ws.Range("D31").Value = "19.83"
$ws.Range("E31").Value = "  -1.21%  "

$ws.Range("D32").Value = "153.31"
$ws.Range("E32").Value = "  -4.48%  "

$ws.Range("D33").Value = "0.0801"
$ws.Range("E33").Value = "  -4.98%  "

$ws.Range("E34").Value = "  +3.55%  "

$ws.Range("E35").Value = "  -3.65%  "

$ws.Range("E36").Value = "  -0.19%  "

$ws.Range("E37").Value = "  -6.56%  "

$ws.Range("E38").Value = "  -9.98%  "

$ws.Range("D39").Value = "14.88"
$ws.Range("E39").Value = "  -10.01%  "

$ws.Range("E40").Value = "  -8.60%  "

$ws.Range("D41").Value = "3.80"
$ws.Range("E41").Value = "  -8.17%  "

$ws.Range("D42").Value = "0.0301"
$ws.Range("E42").Value = "  -3.95%  "

$ws.Range("E43").Value = "  +0.39%  "

$ws.Range("D44").Value = "1.733.58"
$ws.Range("E44").Value = "  -2.19%  "

$ws.Range("D45").Value = "86.34"
$ws.Range("E45").Value = "  +6.58%  "

$ws.Range("D46").Value = "0.187"
$ws.Range("E46").Value = "  -4.15%  "

$ws.Range("D47").Value = "99.78"
$ws.Range("E47").Value = "  -3.77%  "

$ws.Range("D48").Value = "4.90"
$ws.Range("E48").Value = "  -5.15%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "14.68"
$ws.Range("E49").Value = "  -0.22%  "

$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").Value = "8.08"
$ws.Range("E50").Value = "  -2.77%  "

$ws.Range("B51").Value = "ordi"
$ws.Range("C51").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D51").Value = "68.78"
$ws.Range("E51").Value = "  -7.90%  "

# Restore default (General) style now that text values are locked in as strings
$colD.Style = "Normal"
